$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -325
$ws.Range("N19").ClearContents()

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 304.5862
$ws.Range("I33").Value = 267.8889
$ws.Range("K33").Value = 267.8889
$ws.Range("M33").Value = -38.88889999999998

# ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 310.27777
$ws.Range("I41").Value = 90
$ws.Range("K41").Value = 90
$ws.Range("M41").Value = 350

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 799.75
$ws.Range("J43").Value = 799.75
$ws.Range("L43").Value = 799.75
$ws.Range("N43").Value = -937.75

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 895
$ws.Range("I70").Value = 894.44446
$ws.Range("J70").Value = 900
$ws.Range("K70").Value = 2683.33338
$ws.Range("L70").Value = 2700
$ws.Range("M70").Value = -2413.33338
$ws.Range("N70").Value = -3240

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 895
$ws.Range("I73").Value = 894.44446
$ws.Range("J73").Value = 900
$ws.Range("K73").Value = 2683.33338
$ws.Range("L73").Value = 2700
$ws.Range("M73").Value = -1747.33338
$ws.Range("N73").Value = -4572

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33120.47
$ws.Range("I32").Value = 37526.965
$ws.Range("J32").Value = 7562.8
$ws.Range("K32").Value = 37526.965
$ws.Range("L32").Value = 7562.8
$ws.Range("M32").Value = -37239.965
$ws.Range("N32").Value = -8136.8

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15539.281
$ws.Range("I31").Value = 32856.5
$ws.Range("J31").Value = 5148.95
$ws.Range("K31").Value = 32856.5
$ws.Range("L31").Value = 5148.95
$ws.Range("M31").Value = -32561.5
$ws.Range("N31").Value = -5738.95

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 15539.281
$ws.Range("I34").Value = 32856.5
$ws.Range("J34").Value = 5148.95
$ws.Range("K34").Value = 32856.5
$ws.Range("L34").Value = 5148.95
$ws.Range("M34").Value = -32654.5
$ws.Range("N34").Value = -5552.95

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 11000
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -1572
$ws.Range("N41").Value = -20856

# CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 11686.692
$ws.Range("J60").Value = 11686.692
$ws.Range("L60").Value = 11686.692
$ws.Range("N60").Value = -12708.692

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 31551.182
$ws.Range("J74").Value = 31551.182
$ws.Range("L74").Value = 31551.182
$ws.Range("N74").Value = -33299.182

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 31551.182
$ws.Range("J77").Value = 31551.182
$ws.Range("L77").Value = 94653.546
$ws.Range("N77").Value = -103389.546

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 977.7222
$ws.Range("I122").Value = 1055
$ws.Range("J122").Value = 900.44446
$ws.Range("K122").Value = 3165
$ws.Range("L122").Value = 2701.33338
$ws.Range("M122").Value = -715
$ws.Range("N122").Value = -7601.33338

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 490
$ws.Range("I4").Value = 490
$ws.Range("K4").Value = 1470
$ws.Range("M4").Value = -1358

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2049.5
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 3300
$ws.Range("L22").Value = 8997
$ws.Range("M22").Value = -3131
$ws.Range("N22").Value = -9335

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 2049.5
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 2999
$ws.Range("K27").Value = 3300
$ws.Range("L27").Value = 8997
$ws.Range("M27").Value = -3198
$ws.Range("N27").Value = -9201

# CUL row 112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 1983.3334
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 2500
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 7500
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -9716

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 16101.429
$ws.Range("I113").Value = 33605.332
$ws.Range("K113").Value = 100815.996
$ws.Range("M113").Value = -98645.99600000001

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 750.87
$ws.Range("J131").Value = 750.87
$ws.Range("L131").Value = 2252.61
$ws.Range("N131").Value = -12332.61

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2395.2307
$ws.Range("J139").Value = 6495
$ws.Range("L139").Value = 19485
$ws.Range("N139").Value = -29765

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20836630
$ws.Range("I70").Value = 4888
$ws.Range("J70").Value = 31252500
$ws.Range("K70").Value = 4888
$ws.Range("L70").Value = 31252500
$ws.Range("M70").Value = -4618
$ws.Range("N70").Value = -31253040

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 20836630
$ws.Range("I73").Value = 4888
$ws.Range("J73").Value = 31252500
$ws.Range("K73").Value = 4888
$ws.Range("L73").Value = 31252500
$ws.Range("M73").Value = -3952
$ws.Range("N73").Value = -31254372

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 83064.42
$ws.Range("I132").Value = 81172
$ws.Range("K132").Value = 243516
$ws.Range("M132").Value = -240986

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1828.3158
$ws.Range("I82").Value = 2149.077
$ws.Range("J82").Value = 1133.3334
$ws.Range("K82").Value = 2149.077
$ws.Range("L82").Value = 1133.3334
$ws.Range("M82").Value = -1788.077
$ws.Range("N82").Value = -1855.3334

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1828.3158
$ws.Range("I85").Value = 2149.077
$ws.Range("J85").Value = 1133.3334
$ws.Range("K85").Value = 2149.077
$ws.Range("L85").Value = 1133.3334
$ws.Range("M85").Value = -901.0770000000002
$ws.Range("N85").Value = -3629.3334

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2171.4285
$ws.Range("I81").Value = 2180
$ws.Range("J81").Value = 2150
$ws.Range("K81").Value = 4360
$ws.Range("L81").Value = 4300
$ws.Range("M81").Value = -3299
$ws.Range("N81").Value = -6422

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2171.4285
$ws.Range("I84").Value = 2180
$ws.Range("J84").Value = 2150
$ws.Range("K84").Value = 21800
$ws.Range("L84").Value = 21500
$ws.Range("M84").Value = -16496
$ws.Range("N84").Value = -32108

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 865.3333
$ws.Range("I100").Value = 447.75
$ws.Range("K100").Value = 895.5
$ws.Range("M100").Value = -354.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2989.875
$ws.Range("I132").Value = 2759.9
$ws.Range("J132").Value = 3373.1667
$ws.Range("K132").Value = 8279.700000000001
$ws.Range("L132").Value = 10119.5001
$ws.Range("M132").Value = -5749.700000000001
$ws.Range("N132").Value = -15179.5001
